$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Populate the header row on Sheet2 with the new merge-field placeholders.
$ws2.Range("A1").Value = '${username}'
$ws2.Range("B1").Value = '${password}'
$ws2.Range("C1").Value = '${membership}'
$ws2.Range("D1").Value = '${subscription_paid_by}'
$ws2.Range("E1").Value = '${subscription_amount}'
$ws2.Range("F1").Value = '${currency}'
$ws2.Range("G1").Value = '${commence_date}'

# Match the column widths from the target workbook.
$ws2.Columns.Item(1).ColumnWidth = 13.140625
$ws2.Columns.Item(2).ColumnWidth = 15.42578125
$ws2.Columns.Item(3).ColumnWidth = 15.5703125
$ws2.Columns.Item(4).ColumnWidth = 22.140625
$ws2.Columns.Item(5).ColumnWidth = 22.42578125
$ws2.Columns.Item(6).ColumnWidth = 11.85546875
$ws2.Columns.Item(7).ColumnWidth = 18

# Sheet2 becomes the active/visible sheet with a scrolled viewport and a new selection.
$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws2.Range("G2").Select()
